# Updates the "Arveja Verde" sheet with a new weekly price record.
# Rows 79-111 each shift down by one record (row N takes the values that
# used to live in row N-1), a brand-new record is written into row 79,
# and the old row 111 record is appended as the new row 112.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D = Fecha, H = Variedad, J = Volumen, K = Precio minimo,
# L = Precio maximo, M = Precio promedio ponderado, O = Origen, P = Precio $/Kg
$records = @(
    @{ Row = 79;  D = 44876; H = "Sin especificar"; J = 500; K = 15000; L = 15000; M = 15000; O = "Región del Maule";          P = 600  },
    @{ Row = 80;  D = 44504; H = "Sin especificar"; J = 600; K = 12000; L = 12000; M = 12000; O = "Región del Maule";          P = 480  },
    @{ Row = 81;  D = 44250; H = "Sin especificar"; J = 200; K = 22000; L = 22000; M = 22000; O = "Región de La Araucanía";    P = 880  },
    @{ Row = 82;  D = 44166; H = "Sin especificar"; J = 400; K = 16000; L = 16000; M = 16000; O = "Región del Maule";          P = 640  },
    @{ Row = 83;  D = 44855; H = "Sin especificar"; J = 150; K = 13000; L = 13000; M = 13000; O = "Región del Maule";          P = 520  },
    @{ Row = 84;  D = 44494; H = "Sin especificar"; J = 300; K = 17000; L = 17000; M = 17000; O = "Región del Maule";          P = 680  },
    @{ Row = 85;  D = 44539; H = "Sin especificar"; J = 200; K = 18000; L = 18000; M = 18000; O = "Región de La Araucanía";    P = 720  },
    @{ Row = 86;  D = 44518; H = "Sin especificar"; J = 600; K = 15000; L = 15000; M = 15000; O = "Región del Maule";          P = 600  },
    @{ Row = 87;  D = 44488; H = "Sin especificar"; J = 200; K = 20000; L = 20000; M = 20000; O = "Región de O'Higgins";       P = 800  },
    @{ Row = 88;  D = 44488; H = "Sin especificar"; J = 200; K = 19000; L = 19000; M = 19000; O = "Región del Maule";          P = 760  },
    @{ Row = 89;  D = 44614; H = "Sin especificar"; J = 200; K = 27000; L = 27000; M = 27000; O = "Carahue";                   P = 1080 },
    @{ Row = 90;  D = 44179; H = "Sin especificar"; J = 300; K = 22000; L = 22000; M = 22000; O = "Región de La Araucanía";    P = 880  },
    @{ Row = 91;  D = 44168; H = "Sin especificar"; J = 300; K = 20000; L = 20000; M = 20000; O = "Región de La Araucanía";    P = 800  },
    @{ Row = 92;  D = 44169; H = "Perfection";      J = 300; K = 20000; L = 20000; M = 20000; O = "Región de La Araucanía";    P = 800  },
    @{ Row = 93;  D = 44260; H = "Sin especificar"; J = 100; K = 22000; L = 22000; M = 22000; O = "Región de La Araucanía";    P = 880  },
    @{ Row = 94;  D = 44517; H = "Perfection";      J = 500; K = 15000; L = 15000; M = 15000; O = "Región del Maule";          P = 600  },
    @{ Row = 95;  D = 44487; H = "Sin especificar"; J = 200; K = 20000; L = 20000; M = 20000; O = "Región de O'Higgins";       P = 800  },
    @{ Row = 96;  D = 44487; H = "Sin especificar"; J = 150; K = 19000; L = 19000; M = 19000; O = "Región del Maule";          P = 760  },
    @{ Row = 97;  D = 44641; H = "Sin especificar"; J = 300; K = 23000; L = 23000; M = 23000; O = "Carahue";                   P = 920  },
    @{ Row = 98;  D = 44167; H = "Sin especificar"; J = 400; K = 16000; L = 16000; M = 16000; O = "Región del Maule";          P = 640  },
    @{ Row = 99;  D = 44498; H = "Sin especificar"; J = 400; K = 15000; L = 15000; M = 15000; O = "Región del Maule";          P = 600  },
    @{ Row = 100; D = 44858; H = "Sin especificar"; J = 150; K = 15000; L = 15000; M = 15000; O = "Región del Maule";          P = 600  },
    @{ Row = 101; D = 44263; H = "Sin especificar"; J = 300; K = 22000; L = 22000; M = 22000; O = "Región de La Araucanía";    P = 880  },
    @{ Row = 102; D = 44484; H = "Sin especificar"; J = 100; K = 22000; L = 22000; M = 22000; O = "Región de O'Higgins";       P = 880  },
    @{ Row = 103; D = 44495; H = "Sin especificar"; J = 300; K = 17000; L = 17000; M = 17000; O = "Región del Maule";          P = 680  },
    @{ Row = 104; D = 44258; H = "Sin especificar"; J = 200; K = 22000; L = 22000; M = 22000; O = "Región de La Araucanía";    P = 880  },
    @{ Row = 105; D = 44603; H = "Sin especificar"; J = 200; K = 28000; L = 28000; M = 28000; O = "Carahue";                   P = 1120 },
    @{ Row = 106; D = 44642; H = "Sin especificar"; J = 200; K = 25000; L = 25000; M = 25000; O = "Carahue";                   P = 1000 },
    @{ Row = 107; D = 44503; H = "Sin especificar"; J = 600; K = 10000; L = 12000; M = 11000; O = "Región del Maule";          P = 440  },
    @{ Row = 108; D = 44490; H = "Sin especificar"; J = 300; K = 20000; L = 20000; M = 20000; O = "Región del Maule";          P = 800  },
    @{ Row = 109; D = 44845; H = "Perfection";      J = 150; K = 28000; L = 28000; M = 28000; O = "Región de O'Higgins";       P = 1120 },
    @{ Row = 110; D = 44497; H = "Sin especificar"; J = 400; K = 16000; L = 16000; M = 16000; O = "Región de O'Higgins";       P = 640  },
    @{ Row = 111; D = 44489; H = "Sin especificar"; J = 200; K = 18000; L = 18000; M = 18000; O = "Región del Maule";          P = 720  },
    @{ Row = 112; D = 44554; H = "Sin especificar"; J = 300; K = 16000; L = 16000; M = 16000; O = "Carahue";                   P = 640  }
)

# Column D (Fecha) uses a custom date-time display format; grab it once from
# an existing date cell so the new row 112 keeps the same formatting.
$dateFormat = $ws.Range("D111").NumberFormat

foreach ($rec in $records) {
    $r = $rec.Row

    $ws.Cells.Item($r, 4).Value = $rec.D   # D - Fecha
    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat

    $ws.Cells.Item($r, 8).Value  = $rec.H  # H - Variedad
    $ws.Cells.Item($r, 10).Value = $rec.J  # J - Volumen
    $ws.Cells.Item($r, 11).Value = $rec.K  # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $rec.L  # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $rec.M  # M - Precio promedio ponderado
    $ws.Cells.Item($r, 15).Value = $rec.O  # O - Origen
    $ws.Cells.Item($r, 16).Value = $rec.P  # P - Precio $/Kg
}

# Row 112 is brand new -- populate the columns that are constant across every
# other data row in this block (A,B,C,E,F,G,I,N,Q,R), matching row 111's
# original (pre-shift) content exactly.
$ws.Cells.Item(112, 1).Value  = 5                                  # A - Mercado ID
$ws.Cells.Item(112, 2).Value  = "Macroferia Regional de Talca"     # B - Mercado
$ws.Cells.Item(112, 3).Value  = "Maule"                            # C - Región
$ws.Cells.Item(112, 5).Value  = 7                                  # E - Codreg
$ws.Cells.Item(112, 6).Value  = 100112022                          # F - Categoría ID
$ws.Cells.Item(112, 7).Value  = "Arveja Verde"                     # G - Categoría
$ws.Cells.Item(112, 9).Value  = "Primera"                          # I - Calidad
$ws.Cells.Item(112, 14).Value = "`$/saco 25 kilos"                  # N - Unidad de comercialización
$ws.Cells.Item(112, 17).Value = 25                                 # Q - Kg o Unidades
$ws.Cells.Item(112, 18).Value = "Hortaliza"                        # R - Clasificación
